$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# Values in column E are percentage strings with surrounding spaces, stored as text.
# Some column D values are plain-looking numbers (e.g. "0.650", "18.40") which Excel's
# Value setter would silently coerce to numeric and normalize (dropping trailing zeros).
# To preserve them as text exactly as in the source data, we temporarily force the cell
# to Text format, assign the value, then restore the default "Normal" style so no stray
# formatting is left behind on the cell.

$ws.Range("D2").Value = '60.801.39'
$ws.Range("E2").Value = '  -1.58%  '
$ws.Range("D3").Value = '2.906.52'
$ws.Range("E3").Value = '  -2.50%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.45%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -3.89%  '
$ws.Range("D9").Value = '2.911.79'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.05'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("D13").Value = '3.410.99'
$ws.Range("E13").Value = '  -2.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.128'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.20%  '
$ws.Range("D15").Value = '60.741.99'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.87%  '
$ws.Range("D17").Value = '2.914.59'
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("E18").Value = '  -4.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.29%  '
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  +1.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("E26").Value = '  -4.31%  '
$ws.Range("E27").Value = '  -6.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.23%  '
$ws.Range("E30").Value = '  -9.03%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("E33").Value = '  -4.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '153.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.37'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.32%  '
$ws.Range("E36").Value = '  -6.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.997'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.07%  '
$ws.Range("E38").Value = '  -5.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.44%  '
$ws.Range("D42").Value = '2.293.39'
$ws.Range("E42").Value = '  -5.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.650'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0582'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.55%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.97'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.17%  '
$ws.Range("E48").Value = '  -3.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("E50").Value = '  -3.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.13%  '
